$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the video file name for row 2 (foot keypoint entry had no value, so collapse to single row)
$ws.Range("B2").Value = "M-situps.mp4"

# Remove rows 3 and 4 entirely, leaving only the header row and the single data row
$ws.Rows("3:4").Delete()
